$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New line-loss values per row (row number => values for columns C,D,E,F,G,I,J,M,N,O)
$newValues = @{
    2 = @(0.008488266239329079, 0.147827731758789, 0.1135573271574626, 1.400038456624614, 0.002424900010206532, 0.7882384439840138, 0.1163116416311372, 1.206544292833215, 1.784330733312515, 3.624970334784166)
    3 = @(0.007891070831320235, 0.1491069583724691, 0.1153581644060511, 1.378893910803924, 0.002429126498985901, 0.7670961059577763, 0.1187886322385552, 1.09172960567868, 1.642726360592775, 3.540307590685927)
    4 = @(0.007521218297984689, 0.1499363750038172, 0.1165243555414466, 1.366889109471728, 0.002431861185832088, 0.7546733297920909, 0.1203946575050487, 1.021101487973297, 1.555982062698575, 3.490890437310213)
    5 = @(0.00736970608735632, 0.1502854356882051, 0.1170147812078239, 1.362242167891523, 0.002433010811566632, 0.7497508901908816, 0.1210704882305862, 0.9922894490069609, 1.520687862725822, 3.471395350939758)
    6 = @(0.007344499725821407, 0.1503440656711272, 0.117097133295389, 1.361485327040754, 0.002433203836562355, 0.7489419590860322, 0.1211839982562353, 0.9875034644110485, 1.514830721491961, 3.468196958834596)
    7 = @(0.007519178163288842, 0.1499410377353017, 0.1165309080932164, 1.366825447836661, 0.002431876547426568, 0.7546063781568151, 0.1204036855631214, 1.020713038885617, 1.555505845302491, 3.490624920344487)
    8 = @(0.008283013487258017, 0.1482596906185996, 0.1141656981385006, 1.392544319202926, 0.002426328396882081, 0.7808322671466499, 0.1171480082433058, 1.166984953600021, 1.735466194695022, 3.595244500033175)
    9 = @(0.009755641532038339, 0.1453108917029016, 0.1100078062522154, 1.450781054322235, 0.00241655101824274, 0.8367227668455257, 0.1114412321123766, 1.452683323597938, 2.089802097659685, 3.820891310451714)
    10 = @(0.01082217687547882, 0.143355955445476, 0.1072463413851232, 1.498388701719477, 0.002410032457149093, 0.8805504301261919, 0.107664179182585, 1.661780514785562, 2.350816333304238, 3.999363216879658)
    11 = @(0.01130403021591064, 0.142512367649321, 0.1060538909735231, 1.521107967837082, 0.002407209822184543, 0.9010992244748763, 0.1060367339144976, 1.756707513715611, 2.46966906419766, 4.083354593130366)
    12 = @(0.01148601625681778, 0.1421994861815872, 0.1056115221743967, 1.529864966532756, 0.00240616136486588, 0.9089691419266899, 0.1054335652627474, 1.792623980930188, 2.51468852519406, 4.115566264745098)
    13 = @(0.01144684378085259, 0.1422665788642348, 0.1057063853499809, 1.527972138919353, 0.002406386262673441, 0.9072702668385659, 0.1055628844276564, 1.784890118116863, 2.504992288223264, 4.10861080482789)
    14 = @(0.01131901198759522, 0.1424864951362286, 0.1060173128123711, 1.521825325449427, 0.002407123156445686, 0.9017449093747985, 0.1059868477920407, 1.759663006085361, 2.473372615325957, 4.085996514550231)
    15 = @(0.01124064851521922, 0.142622055176469, 0.1062089617409128, 1.518080269959796, 0.002407577180612424, 0.8983720167590974, 0.1062482466045118, 1.744206643302647, 2.45400615902696, 4.072197560400582)
    16 = @(0.010790619489363, 0.1434120052682972, 0.1073255552976518, 1.496925401401029, 0.00241021978539305, 0.8792198733739696, 0.1077723678236264, 1.655572713001007, 2.343051004704932, 3.993930837959624)
    17 = @(0.01051368827565113, 0.1439083190415928, 0.1080268913235405, 1.484220328986396, 0.002411877408996095, 0.8676276805007461, 0.1087306532838053, 1.601147495867295, 2.275010344358236, 3.946636799236046)
    18 = @(0.01035409261719877, 0.1441980896884514, 0.1084362809483383, 1.477012639847317, 0.002412844267450758, 0.8610176654635211, 0.1092903724761527, 1.569825602467191, 2.23588632496768, 3.919698239223635)
    19 = @(0.01030000278318255, 0.1442969405230428, 0.1085759231296799, 1.474589374899978, 0.002413173940212523, 0.8587894861112915, 0.1094813485190222, 1.559217546525247, 2.22264165066963, 3.91062251474716)
    20 = @(0.01054320041282608, 0.143855040194514, 0.1079516118286188, 1.485562456189982, 0.002411699562096265, 0.8688557329248425, 0.1086277580094119, 1.606943030065082, 2.282252263305168, 3.951644016177795)
    21 = @(0.0113565723817004, 0.1424217222347988, 0.105925736417535, 1.523626614454457, 0.002406906159684641, 0.9033654330135192, 0.1058619632128321, 1.767073664613548, 2.482659773179762, 4.092627839675799)
    22 = @(0.01188535033613647, 0.1415232395141148, 0.1046552631894868, 1.549400101355033, 0.002403892330147814, 0.9264358581785501, 0.1041307966227327, 1.871550299196741, 2.613708692680348, 4.187136482930725)
    23 = @(0.01160338984271903, 0.1419992772319489, 0.1053284323419059, 1.53556197464647, 0.00240549001961083, 0.9140752865017276, 0.1050477389061326, 1.815806314762682, 2.543760284409359, 4.136477876668323)
    24 = @(0.01052985917290528, 0.1438791137516173, 0.1079856264326651, 1.484955379571375, 0.002411779923527806, 0.8683003607368818, 0.1086742495931841, 1.604322968870179, 2.278978211791753, 3.94937947048345)
    25 = @(0.009359961063456268, 0.1460714127853679, 0.1110811530826048, 1.434184538676519, 0.002419078777137697, 0.8211206256214325, 0.1129122541046845, 1.375527898481025, 1.993811236618569, 3.757633974406986)
}

$columns = @("C", "D", "E", "F", "G", "I", "J", "M", "N", "O")

foreach ($row in $newValues.Keys) {
    $values = $newValues[$row]
    for ($i = 0; $i -lt $columns.Count; $i++) {
        $ws.Range($columns[$i] + $row).Value = $values[$i]
    }
}
